$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The data rows now all refer to the "testValidLogin3" test case (TC-TESTVALIDLOGIN3),
# executed multiple times in a group run, alternating SKIP/FAIL status.
$testId = "TC-TESTVALIDLOGIN3"
$methodName = "testValidLogin3"
$statuses = @("SKIP", "FAIL", "SKIP", "FAIL", "SKIP", "FAIL")

for ($i = 0; $i -lt $statuses.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 1).Value = $testId
    $ws.Cells.Item($row, 2).Value = $methodName
    $ws.Cells.Item($row, 3).Value = $statuses[$i]
}
